$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update "想去人数" (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 501
$ws1.Range("F3").Value = 5919
$ws1.Range("F4").Value = 392
$ws1.Range("F5").Value = 79
$ws1.Range("F6").Value = 105

# Sheet "全部类型" (all types): update "想去人数" (F column) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 501
$ws4.Range("F3").Value = 5919
$ws4.Range("F4").Value = 392
$ws4.Range("F6").Value = 79
$ws4.Range("F7").Value = 105
